$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.223446
$ws.Range("H2").Value = 87.670338
$ws.Range("I2").Value = 0.0169041244192178
$ws.Range("J2").Value = 0.0169041244192178
$ws.Range("M2").Value = 2.893186333333333
$ws.Range("N2").Value = 8.679558999999999
$ws.Range("O2").Value = 0.296307560753476
$ws.Range("P2").Value = 0.296307560753476
$ws.Range("Q2").Value = 84.54887458010465
$ws.Range("R2").Value = 760.939871220942
$ws.Range("S2").Value = 0.005008819873331695
$ws.Range("T2").Value = 0.005008819873331695

$ws.Range("G3").Value = 29.223446
$ws.Range("H3").Value = 87.670338
$ws.Range("I3").Value = 0.0169041244192178
$ws.Range("J3").Value = 0.0169041244192178
$ws.Range("O3").Value = 0.2593826903485334
$ws.Range("P3").Value = 0.2593826903485334
$ws.Range("Q3").Value = 74.01267284156199
$ws.Range("R3").Value = 666.1140555740579
$ws.Range("S3").Value = 0.004384637269843054
$ws.Range("T3").Value = 0.004384637269843054

$ws.Range("G4").Value = 29.223446
$ws.Range("H4").Value = 87.670338
$ws.Range("I4").Value = 0.0169041244192178
$ws.Range("J4").Value = 0.0169041244192178
$ws.Range("M4").Value = 3.448217666666667
$ws.Range("N4").Value = 10.344653
$ws.Range("O4").Value = 0.353151455882854
$ws.Range("P4").Value = 0.353151455882854
$ws.Range("Q4").Value = 100.7688027780793
$ws.Range("R4").Value = 906.9192250027141
$ws.Range("S4").Value = 0.005969716149071671
$ws.Range("T4").Value = 0.00596971614907167

$ws.Range("G5").Value = 29.223446
$ws.Range("H5").Value = 87.670338
$ws.Range("I5").Value = 0.0169041244192178
$ws.Range("J5").Value = 0.0169041244192178
$ws.Range("M5").Value = 0.8900816666666667
$ws.Range("N5").Value = 2.670245
$ws.Range("O5").Value = 0.09115829301513656
$ws.Range("P5").Value = 0.09115829301513655
$ws.Range("Q5").Value = 26.01125352142333
$ws.Range("R5").Value = 234.10128169281
$ws.Range("S5").Value = 0.001540951126971382
$ws.Range("T5").Value = 0.001540951126971381

$ws.Range("I6").Value = 0.9471112884046843
$ws.Range("J6").Value = 0.9471112884046842
$ws.Range("M6").Value = 2.893186333333333
$ws.Range("N6").Value = 8.679558999999999
$ws.Range("O6").Value = 0.296307560753476
$ws.Range("P6").Value = 0.296307560753476
$ws.Range("Q6").Value = 4737.139383906307
$ws.Range("R6").Value = 42634.25445515676
$ws.Range("S6").Value = 0.2806362356292739
$ws.Range("T6").Value = 0.2806362356292739

$ws.Range("I7").Value = 0.9471112884046843
$ws.Range("J7").Value = 0.9471112884046842
$ws.Range("O7").Value = 0.2593826903485334
$ws.Range("P7").Value = 0.2593826903485334
$ws.Range("Q7").Value = 4146.812706463136
$ws.Range("S7").Value = 0.2456642740458727
$ws.Range("T7").Value = 0.2456642740458727

$ws.Range("I8").Value = 0.9471112884046843
$ws.Range("J8").Value = 0.9471112884046842
$ws.Range("M8").Value = 3.448217666666667
$ws.Range("N8").Value = 10.344653
$ws.Range("O8").Value = 0.353151455882854
$ws.Range("P8").Value = 0.353151455882854
$ws.Range("Q8").Value = 5645.916242881065
$ws.Range("R8").Value = 50813.24618592959
$ws.Range("S8").Value = 0.3344737303831999
$ws.Range("T8").Value = 0.3344737303831998

$ws.Range("I9").Value = 0.9471112884046843
$ws.Range("J9").Value = 0.9471112884046842
$ws.Range("M9").Value = 0.8900816666666667
$ws.Range("N9").Value = 2.670245
$ws.Range("O9").Value = 0.09115829301513656
$ws.Range("P9").Value = 0.09115829301513655
$ws.Range("Q9").Value = 1457.369291939705
$ws.Range("R9").Value = 13116.32362745735
$ws.Range("S9").Value = 0.08633704834633772
$ws.Range("T9").Value = 0.08633704834633769

$ws.Range("G10").Value = 37.39212666666667
$ws.Range("H10").Value = 112.17638
$ws.Range("I10").Value = 0.02162924801792661
$ws.Range("J10").Value = 0.0216292480179266
$ws.Range("M10").Value = 2.893186333333333
$ws.Range("N10").Value = 8.679558999999999
$ws.Range("O10").Value = 0.296307560753476
$ws.Range("P10").Value = 0.296307560753476
$ws.Range("Q10").Value = 108.1823898462689
$ws.Range("R10").Value = 973.64150861642
$ws.Range("S10").Value = 0.006408909721123788
$ws.Range("T10").Value = 0.006408909721123787

$ws.Range("G11").Value = 37.39212666666667
$ws.Range("H11").Value = 112.17638
$ws.Range("I11").Value = 0.02162924801792661
$ws.Range("J11").Value = 0.0216292480179266
$ws.Range("O11").Value = 0.2593826903485334
$ws.Range("P11").Value = 0.2593826903485334
$ws.Range("Q11").Value = 94.70105742595334
$ws.Range("R11").Value = 852.30951683358
$ws.Range("S11").Value = 0.005610252541105487
$ws.Range("T11").Value = 0.005610252541105486

$ws.Range("G12").Value = 37.39212666666667
$ws.Range("H12").Value = 112.17638
$ws.Range("I12").Value = 0.02162924801792661
$ws.Range("J12").Value = 0.0216292480179266
$ws.Range("M12").Value = 3.448217666666667
$ws.Range("N12").Value = 10.344653
$ws.Range("O12").Value = 0.353151455882854
$ws.Range("P12").Value = 0.353151455882854
$ws.Range("Q12").Value = 128.9361917662378
$ws.Range("R12").Value = 1160.42572589614
$ws.Range("S12").Value = 0.007638400427182116
$ws.Range("T12").Value = 0.007638400427182114

$ws.Range("G13").Value = 37.39212666666667
$ws.Range("H13").Value = 112.17638
$ws.Range("I13").Value = 0.02162924801792661
$ws.Range("J13").Value = 0.0216292480179266
$ws.Range("M13").Value = 0.8900816666666667
$ws.Range("N13").Value = 2.670245
$ws.Range("O13").Value = 0.09115829301513656
$ws.Range("P13").Value = 0.09115829301513655
$ws.Range("Q13").Value = 33.28204642367778
$ws.Range("R13").Value = 299.5384178131
$ws.Range("S13").Value = 0.001971685328515215
$ws.Range("T13").Value = 0.001971685328515215

$ws.Range("G14").Value = 24.817167
$ws.Range("H14").Value = 74.45150100000001
$ws.Range("I14").Value = 0.01435533915817136
$ws.Range("J14").Value = 0.01435533915817136
$ws.Range("M14").Value = 2.893186333333333
$ws.Range("N14").Value = 8.679558999999999
$ws.Range("O14").Value = 0.296307560753476
$ws.Range("P14").Value = 0.296307560753476
$ws.Range("Q14").Value = 71.80068839645099
$ws.Range("R14").Value = 646.206195568059
$ws.Range("S14").Value = 0.004253595529746614
$ws.Range("T14").Value = 0.004253595529746613

$ws.Range("G15").Value = 24.817167
$ws.Range("H15").Value = 74.45150100000001
$ws.Range("I15").Value = 0.01435533915817136
$ws.Range("J15").Value = 0.01435533915817136
$ws.Range("O15").Value = 0.2593826903485334
$ws.Range("P15").Value = 0.2593826903485334
$ws.Range("Q15").Value = 62.853123551049
$ws.Range("R15").Value = 565.678111959441
$ws.Range("S15").Value = 0.003723526491712138
$ws.Range("T15").Value = 0.003723526491712138

$ws.Range("G16").Value = 24.817167
$ws.Range("H16").Value = 74.45150100000001
$ws.Range("I16").Value = 0.01435533915817136
$ws.Range("J16").Value = 0.01435533915817136
$ws.Range("M16").Value = 3.448217666666667
$ws.Range("N16").Value = 10.344653
$ws.Range("O16").Value = 0.353151455882854
$ws.Range("P16").Value = 0.353151455882854
$ws.Range("Q16").Value = 85.57499368601701
$ws.Range("R16").Value = 770.1749431741531
$ws.Range("S16").Value = 0.00506960892340036
$ws.Range("T16").Value = 0.005069608923400359

$ws.Range("G17").Value = 24.817167
$ws.Range("H17").Value = 74.45150100000001
$ws.Range("I17").Value = 0.01435533915817136
$ws.Range("J17").Value = 0.01435533915817136
$ws.Range("M17").Value = 0.8900816666666667
$ws.Range("N17").Value = 2.670245
$ws.Range("O17").Value = 0.09115829301513656
$ws.Range("P17").Value = 0.09115829301513655
$ws.Range("Q17").Value = 22.089305365305
$ws.Range("R17").Value = 198.803748287745
$ws.Range("S17").Value = 0.001308608213312249
$ws.Range("T17").Value = 0.001308608213312248

